# Fruta / hortaliza, semanal
# The weekly refresh re-sorted the detail rows (2-13) of the "Haba" sheet.
# Every column for a given row moved together, so capture the full rows
# first and then write them back out in their new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot the existing rows (2..13) before overwriting anything.
$data = @{}
foreach ($r in 2..13) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$r").Value2
    }
    $data[$r] = $rowData
}

# Destination row -> source row (where the data used to live before the edit).
$mapping = @{
    2  = 12
    3  = 13
    4  = 2
    5  = 11
    6  = 6
    7  = 10
    8  = 9
    9  = 5
    10 = 4
    11 = 3
    12 = 7
    13 = 8
}

foreach ($destRow in 2..13) {
    $srcRow = $mapping[$destRow]
    $srcData = $data[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcData[$col]
    }
}
